# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# listing with freshly scraped values. Numeric-looking price strings are
# forced to Text format before assignment (and restored to the default
# "Normal" style afterwards) so Excel doesn't silently reinterpret them
# as numbers, which would also strip the thousands-dot formatting used
# by some of the other rows (e.g. "70.758.88").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.758.88'
$ws.Range("E2").Value = '  +1.80%  '
$ws.Range("D3").Value = '3.803.04'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  +0.16%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '697.78'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  +7.97%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '172.76'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  +3.97%  '
$ws.Range("D7").Value = '3.802.40'
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("E11").Value = '  +4.68%  '
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("E13").Value = '  +7.76%  '
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '36.18'
$c.Style = 'Normal'
$ws.Range("E14").Value = '  +3.37%  '
$ws.Range("D15").Value = '4.443.01'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '3.791.81'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").Value = '70.806.86'
$ws.Range("E17").Value = '  +2.02%  '
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '17.75'
$c.Style = 'Normal'
$ws.Range("E18").Value = '  -0.22%  '
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '7.18'
$c.Style = 'Normal'
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("E20").Value = '  +0.14%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '11.17'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  +16.36%  '
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '478.29'
$c.Style = 'Normal'
$ws.Range("E22").Value = '  +2.12%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '0.710'
$c.Style = 'Normal'
$ws.Range("E23").Value = '  +0.16%  '
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '83.71'
$c.Style = 'Normal'
$ws.Range("E24").Value = '  +2.24%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '0.0000143'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  -0.98%  '
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '12.32'
$c.Style = 'Normal'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  +2.20%  '
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '10.40'
$c.Style = 'Normal'
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").Value = '3.954.73'
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("E30").Value = '  -0.18%  '
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '3.12'
$c.Style = 'Normal'
$ws.Range("E31").Value = '  +15.07%  '
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("E33").Value = '  +4.79%  '
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '0.188'
$c.Style = 'Normal'
$ws.Range("E34").Value = '  +8.45%  '
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '29.42'
$c.Style = 'Normal'
$ws.Range("E35").Value = '  +2.30%  '
$ws.Range("E36").Value = '  +3.97%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  +2.37%  '
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '3.41'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  +3.50%  '
$ws.Range("E40").Value = '  +2.56%  '
$ws.Range("E41").Value = '  +13.24%  '
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '0.979'
$c.Style = 'Normal'
$ws.Range("E42").Value = '  +2.19%  '
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '0.000327'
$c.Style = 'Normal'
$ws.Range("E43").Value = '  +20.21%  '
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("E45").Value = '  +0.01%  '
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '163.26'
$c.Style = 'Normal'
$ws.Range("E46").Value = '  +3.76%  '
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '48.88'
$c.Style = 'Normal'
$ws.Range("E47").Value = '  +2.62%  '
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '44.38'
$c.Style = 'Normal'
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("E50").Value = '  -1.61%  '
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '8.57'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  +2.05%  '
